$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'SAEP currently in service'
$ws.Range('C2').Value = 'Indicate if the SAEP is currently in service'
$ws.Range('B3').Value = 'SAEP Supply Service According to a Pre-Established Schedule'
$ws.Range('C3').Value = 'The degree to which the SAEP is delivering water as per the predefined timetable to ensure regular and predictable water distribution.'
$ws.Range('B4').Value = 'The SAEP has production meters'
$ws.Range('C4').Value = 'The presence of meters that quantify the total water produced by the SAEP.'
$ws.Range('B5').Value = 'The SAEP has distribution meters'
$ws.Range('C5').Value = 'The presence of meters that measure the volume of water distributed through the network.'
$ws.Range('B6').Value = 'The SAEP has connection/consumption meters'
$ws.Range('C6').Value = 'The availability of meters that record the amount of water consumed by end users.'
$ws.Range('B7').Value = 'SAEP Capacity Used'
$ws.Range('C7').Value = 'The percentage of the total production capacity of the SAEP that is actually being utilized to serve the population.'
$ws.Range('B8').Value = 'SAEP Coverage Rate'
$ws.Range('C8').Value = 'The proportion of the target population that has access to the services provided by the SAEP within its operational area.'
$ws.Range('B9').Value = '# of Employees'
$ws.Range('C9').Value = 'The total number of personnel employed by the SAEP organization.'
$ws.Range('B10').Value = '# of Active Subscribers at the Start of the Month'
$ws.Range('C10').Value = 'The number of active subscribers (those receiving and paying for water services) at the beginning of the month.'
$ws.Range('B11').Value = '# of Active Subscribers at the End of the Month'
$ws.Range('C11').Value = 'The number of active subscribers at the end of the month.'
$ws.Range('B12').Value = '# of Passive Subscribers at the End of the Month'
$ws.Range('C12').Value = 'The number of subscribers who are registered but are not actively receiving water services at the end of the month.'
$ws.Range('B13').Value = '# of Suspended Subscribers at the End of the Month'
$ws.Range('C13').Value = 'The number of subscribers whose services have been temporarily suspended by the end of the month.'
$ws.Range('B14').Value = '# of Subscribers with Arrears'
$ws.Range('C14').Value = 'The number of subscribers who have overdue payments for water services.'
$ws.Range('B15').Value = '# of Subscribers without Water at the End of the Month'
$ws.Range('C15').Value = 'The number of subscribers who did not receive any water supply by the end of the month.'
$ws.Range('B16').Value = 'Total # of Subscribers'
$ws.Range('C16').Value = 'The total count of subscribers, including active, passive, and suspended.'
$ws.Range('B17').Value = 'Residentials'
$ws.Range('C17').Value = 'The number of subscribers classified as residential households.'
$ws.Range('B18').Value = 'Institutionals'
$ws.Range('C18').Value = 'The number of subscribers classified as institutions (schools, hospitals, etc.).'
$ws.Range('B19').Value = 'Residential Flat Rate'
$ws.Range('C19').Value = 'The flat rate charged to residential subscribers for water services.'
$ws.Range('B20').Value = 'Commercial Flat Rate'
$ws.Range('C20').Value = 'The flat rate charged to commercial subscribers for water services.'
$ws.Range('B21').Value = 'Institutional Flat Rate'
$ws.Range('C21').Value = 'The flat rate charged to institutional subscribers for water services.'
$ws.Range('B22').Value = 'Commercials'
$ws.Range('C22').Value = 'The number of subscribers classified as commercial businesses.'
$ws.Range('B23').Value = '# of Meters Disconnected'
$ws.Range('C23').Value = 'The number of water meters that have been disconnected either due to non-payment, malfunction, or other reasons.'
$ws.Range('B24').Value = '# of Meters Reconnected'
$ws.Range('C24').Value = 'The number of water meters that have been reconnected after being temporarily disconnected.'
$ws.Range('B25').Value = '# of illegal connections'
$ws.Range('C25').Value = 'The number of instances where unauthorized water connections have been detected.'
$ws.Range('B26').Value = '# of illegal connections regularized'
$ws.Range('C26').Value = 'The number of previously unauthorized water intakes that have been regularized (made legal).'
$ws.Range('B27').Value = '# of Clients with Active Meters'
$ws.Range('C27').Value = 'The number of subscribers who have meters that are currently operating and providing data.'
$ws.Range('B28').Value = '# of Existing Meters'
$ws.Range('C28').Value = 'The total number of water meters currently installed within the SAEP''s service area.'
$ws.Range('B29').Value = '# of New Installed Meters'
$ws.Range('C29').Value = 'The number of new water meters installed during a specific period.'
$ws.Range('B30').Value = '# of Meters Down'
$ws.Range('C30').Value = 'The number of water meters that are malfunctioning or not operational.'
$ws.Range('B31').Value = '# of Meters Replaced'
$ws.Range('C31').Value = 'The number of old or faulty water meters that have been replaced with new ones.'
$ws.Range('B32').Value = '# of Meters Repaired'
$ws.Range('C32').Value = 'The number of defective water meters that have been repaired and made operational.'
$ws.Range('B33').Value = '# of Clients with Meters'
$ws.Range('C33').Value = 'The number of subscribers who have water meters installed at their premises.'
$ws.Range('B34').Value = '# of Households Served by the SAEP'
$ws.Range('C34').Value = 'The total number of households receiving water services from the SAEP.'
$ws.Range('B35').Value = '# of Individuals Served by the SAEP'
$ws.Range('C35').Value = 'The total number of individuals benefiting from the water services provided by the SAEP.'
$ws.Range('B36').Value = '# of SAEP Water Production per Household (liters/household/month)'
$ws.Range('C36').Value = 'The average quantity of water produced by the SAEP per household per month, measured in liters.'
$ws.Range('B37').Value = 'Total Water Produced (m³/month)'
$ws.Range('C37').Value = 'The total volume of water produced by the SAEP in a month, measured in cubic meters.'
$ws.Range('B38').Value = 'Total Water Distributed (m³/month)'
$ws.Range('C38').Value = 'The total volume of water distributed through the SAEP''s network in a month, measured in cubic meters.'
$ws.Range('B39').Value = 'Total Water Consumed (m³/month)'
$ws.Range('C39').Value = 'The total volume of water consumed by end users (households, businesses, institutions) in a month, measured in cubic meters.'
$ws.Range('B40').Value = '% of Water Loss'
$ws.Range('C40').Value = 'The percentage of water lost (due to leaks, theft, or unaccounted-for usage) compared to the total water produced. Calculated as [(Total Water Produced - Total Water Distributed) / Total Water Produced] * 100.'
$ws.Range('B41').Value = '% of E. Coli Tests Conform'
$ws.Range('C41').Value = 'The percentage of water quality tests for E. Coli that meet the required safety standards.'
$ws.Range('B42').Value = '% of Residual Chlorine Tests Conforming to Norms'
$ws.Range('C42').Value = 'The percentage of water quality tests for residual chlorine that comply with established safety norms.'
$ws.Range('B43').Value = 'Number of Residual Chlorine Tests Conducted'
$ws.Range('C43').Value = 'The total number of residual chlorine tests conducted to monitor water quality.'
$ws.Range('B44').Value = 'Total Revenue (HTG)'
$ws.Range('C44').Value = 'The total income generated by the SAEP from all sources, measured in Haitian Gourdes (HTG).'
$ws.Range('B45').Value = 'Total Expenses (HTG)'
$ws.Range('C45').Value = 'The total operating expenses incurred by the SAEP, measured in Haitian Gourdes (HTG).'
$ws.Range('B46').Value = 'Operating Ratio [(Expenses/Revenue) * 100]'
$ws.Range('C46').Value = 'The operating ratio for a water system is a financial performance metric that indicates how well the system is managing its operating expenses relative to its operating revenues. If the operating ratio is less than 100%: The system is generating more revenue than its operating costs, which indicates financial health and efficiency. If the operating ratio is equal to 100%: The system is breaking even, meaning its revenues are exactly covering its operating expenses. If the operating ratio is greater than 100%: The system is spending more than it earns in revenue, signaling potential financial sustainability issues.'
$ws.Range('B47').Value = 'Amount Collected from Subscribers (HTG)'
$ws.Range('C47').Value = 'The total amount of money collected from subscribers for water services during a specific period, measured in Haitian Gourdes (HTG).'
$ws.Range('B48').Value = 'Amount Billed (HTG)'
$ws.Range('C48').Value = 'The total value of bills issued to subscribers for water services during a specific period, measured in Haitian Gourdes (HTG).'
$ws.Range('B49').Value = 'Collection Efficiency [(Amount Collected / Amount Billed) * 100]'
$ws.Range('C49').Value = 'The efficiency of the SAEP in collecting billed amounts from subscribers. Calculated as (Amount Collected / Amount Billed) * 100.'
$ws.Range('B50').Value = 'Amount Collected Arrears'
$ws.Range('C50').Value = 'The total amount of overdue payments collected from subscribers during a specific period, measured in Haitian Gourdes (HTG).'
$ws.Range('B51').Value = 'Amount Arrears (Existing at the Start of the Month)'
$ws.Range('C51').Value = 'The total amount of overdue payments owed by subscribers at the beginning of the month, measured in Haitian Gourdes (HTG).'
$ws.Range('B52').Value = 'Collection Efficiency - Arrears [(Amount Arrears Collected / Amount Arrears at the Start of the Month) * 100]'
$ws.Range('C52').Value = 'The efficiency of the SAEP in collecting overdue payments from subscribers. Calculated as (Amount Arrears Collected / Amount Arrears at the Start of the Month) * 100.'
$ws.Range('A53:C56').EntireRow.Delete()
